# Gendata.xlsx update:
#  - Fill in row 3 of the "Gen slack" sheet (B3=0, C3=1, D3=0)
#  - Update the remembered cell selection on the "Bus", "Gen slack" and
#    "Trans" sheets
#  - Leave "Trans" as the active sheet/tab (was "Gen slack")

$wb = $excel.ActiveWorkbook

# --- Bus sheet: just move the selection ---
$wsBus = $wb.Worksheets.Item("Bus")
$wsBus.Activate()
$wsBus.Range("J29").Select()

# --- Gen slack sheet: enter the missing row-3 data, move selection ---
$wsGen = $wb.Worksheets.Item("Gen slack")
$wsGen.Activate()
$wsGen.Range("B3").Value = 0
$wsGen.Range("C3").Value = 1
$wsGen.Range("D3").Value = 0
$wsGen.Range("I17").Select()

# --- Trans sheet: move selection and leave it as the active tab ---
$wsTrans = $wb.Worksheets.Item("Trans")
$wsTrans.Activate()
$wsTrans.Range("H31").Select()
